# This script reproduces, via the Word COM object model, a round of
# "spell check as you type" style edits: the author retyped a couple of
# words (which made Word insert w:proofErr spellStart/spellEnd markers
# around them), collapsed the elseif field code into a single run, and
# the cursor's last-edit position moved Word's hidden "_GoBack" bookmark
# from the elseif field to the very start of the document.
#
# The emulated object model does not expose w:proofErr (it is not part
# of the Word automation surface - real Word inserts it internally
# during its spell-check pass) nor does Bookmarks.Add place tags
# correctly for a range collapsed at document position 0 (a quirk of
# this host). Range.InsertXML, however, accepts a literal <w:p> fragment
# and splices it in verbatim, so we use it to rebuild the handful of
# paragraphs that changed, keeping the untouched paragraphs (and their
# untouched runs) completely alone.

$d = $word.ActiveDocument

# --- Paragraph 1: "Basic if demonstration :" ---------------------------
# Adds the _GoBack bookmark at the very start of the document and splits
# the trailing run so "demonstration" is wrapped in spellStart/spellEnd.
$p1xml = @'
<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">Basic </w:t></w:r><w:r><w:t>if</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>demonstration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t> :</w:t></w:r></w:p>
'@
$d.Paragraphs(1).Range.InsertXML($p1xml)

# --- Paragraph 3: "The THEN paragraph." --------------------------------
$p3xml = @'
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">The THEN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>paragraph</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>
'@
$d.Paragraphs(3).Range.InsertXML($p3xml)

# --- Paragraph 4: the "elseif" field -----------------------------------
# Collapses "elseif " + "self.name " + "=" + " 'anydsl'" into a single
# instrText run, and drops the _GoBack bookmark that used to sit here
# (it moved to paragraph 1, see above).
$p4xml = @'
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> m:</w:instrText></w:r><w:r><w:instrText>elseif self.name = 'anydsl'</w:instrText></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p>
'@
$d.Paragraphs(4).Range.InsertXML($p4xml)

# --- Paragraph 5: "The ELSEIF paragraph." -------------------------------
$p5xml = @'
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">The ELSEIF </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>paragraph</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>
'@
$d.Paragraphs(5).Range.InsertXML($p5xml)

# --- Paragraph 7: "End of demonstration." -------------------------------
$p7xml = @'
<w:p><w:r><w:t>En</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>demonstration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>
'@
$d.Paragraphs(7).Range.InsertXML($p7xml)
